$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-12-10T19:06:45"

$ws.Range("W4").Value = 177.02
$ws.Range("X4").Value = 137.5
$ws.Range("Z4").Value = 87.84999999999999
$ws.Range("W6").Value = -9.380000000000001
$ws.Range("X6").Value = -7.01
$ws.Range("Z6").Value = -2.55
$ws.Range("W9").Value = 183.11
$ws.Range("X9").Value = 144.37
$ws.Range("Y9").Value = 119.94
$ws.Range("Z9").Value = 93.87
$ws.Range("W11").Value = -3.3
$ws.Range("X11").Value = -0.14
$ws.Range("Y11").Value = 2.76
$ws.Range("Z11").Value = 3.47
$ws.Range("W14").Value = 183.29
$ws.Range("X14").Value = 144.37
$ws.Range("Y14").Value = 119.94
$ws.Range("Z14").Value = 93.97
$ws.Range("W16").Value = -3.12
$ws.Range("X16").Value = -0.14
$ws.Range("Y16").Value = 2.76
$ws.Range("Z16").Value = 3.57
$ws.Range("W19").Value = 62.41
$ws.Range("X19").Value = 136.98
$ws.Range("Z19").Value = 87.84999999999999
$ws.Range("W20").Value = -113.77
$ws.Range("X20").Value = 0
$ws.Range("W21").Value = -10.22
$ws.Range("X21").Value = -7.53
$ws.Range("Z21").Value = -2.55
$ws.Range("W24").Value = 176.18
$ws.Range("X24").Value = 136.98
$ws.Range("Z24").Value = 87.84999999999999
$ws.Range("W26").Value = -10.22
$ws.Range("X26").Value = -7.53
$ws.Range("Z26").Value = -2.55
$ws.Range("W29").Value = 60.76
$ws.Range("X29").Value = 135.69
$ws.Range("Z29").Value = 87.34
$ws.Range("W30").Value = -113.77
$ws.Range("X30").Value = 0
$ws.Range("W31").Value = -11.87
$ws.Range("X31").Value = -8.82
$ws.Range("Z31").Value = -3.06
$ws.Range("W34").Value = 189.62
$ws.Range("X34").Value = 150.85
$ws.Range("Y34").Value = 126.54
$ws.Range("Z34").Value = 99.01000000000001
$ws.Range("W36").Value = 3.22
$ws.Range("X36").Value = 6.34
$ws.Range("Y36").Value = 9.359999999999999
$ws.Range("Z36").Value = 8.609999999999999
$ws.Range("W39").Value = 177.02
$ws.Range("X39").Value = 137.5
$ws.Range("Z39").Value = 87.84999999999999
$ws.Range("W41").Value = -9.380000000000001
$ws.Range("X41").Value = -7.01
$ws.Range("Z41").Value = -2.55
$ws.Range("W44").Value = 182.93
$ws.Range("X44").Value = 141.68
$ws.Range("Y44").Value = 115.11
$ws.Range("Z44").Value = 89.15000000000001
$ws.Range("W46").Value = -3.48
$ws.Range("X46").Value = -2.83
$ws.Range("Y46").Value = -2.07
$ws.Range("Z46").Value = -1.25
$ws.Range("W49").Value = 199.36
$ws.Range("X49").Value = 154.72
$ws.Range("Y49").Value = 124.52
$ws.Range("Z49").Value = 96.89
$ws.Range("W51").Value = 12.96
$ws.Range("X51").Value = 10.21
$ws.Range("Y51").Value = 7.35
$ws.Range("Z51").Value = 6.49
$ws.Range("W54").Value = 191.97
$ws.Range("X54").Value = 149.91
$ws.Range("Y54").Value = 122.57
$ws.Range("Z54").Value = 95.16
$ws.Range("W56").Value = 5.57
$ws.Range("X56").Value = 5.4
$ws.Range("Y56").Value = 5.39
$ws.Range("Z56").Value = 4.76
$ws.Range("W59").Value = 194.57
$ws.Range("X59").Value = 150.85
$ws.Range("Z59").Value = 94.17
$ws.Range("W61").Value = 8.17
$ws.Range("X61").Value = 6.34
$ws.Range("Z61").Value = 3.77
$ws.Range("W64").Value = 198.51
$ws.Range("X64").Value = 154.06
$ws.Range("Z64").Value = 95.97
$ws.Range("W66").Value = 12.11
$ws.Range("X66").Value = 9.550000000000001
$ws.Range("Z66").Value = 5.57
$ws.Range("X69").Value = 155.39
$ws.Range("Z69").Value = 96.68000000000001
$ws.Range("W71").Value = 13.6
$ws.Range("X71").Value = 10.88
$ws.Range("Z71").Value = 6.28
$ws.Range("W74").Value = 196.63
$ws.Range("X74").Value = 152.76
$ws.Range("Z74").Value = 95.16
$ws.Range("W76").Value = 10.22
$ws.Range("X76").Value = 8.25
$ws.Range("Z76").Value = 4.76
$ws.Range("W79").Value = 196.84
$ws.Range("X79").Value = 153.19
$ws.Range("Z79").Value = 95.43000000000001
$ws.Range("W81").Value = 10.44
$ws.Range("X81").Value = 8.68
$ws.Range("Z81").Value = 5.03
$ws.Range("W84").Value = 189.24
$ws.Range("X84").Value = 147.91
$ws.Range("Y84").Value = 120.68
$ws.Range("Z84").Value = 94.06999999999999
$ws.Range("W86").Value = 2.84
$ws.Range("X86").Value = 3.4
$ws.Range("Y86").Value = 3.5
$ws.Range("Z86").Value = 3.67
$ws.Range("W89").Value = 174.53
$ws.Range("X89").Value = 135.69
$ws.Range("Z89").Value = 87.34
$ws.Range("W91").Value = -11.87
$ws.Range("X91").Value = -8.82
$ws.Range("Z91").Value = -3.06
